$wb = $excel.ActiveWorkbook

# Rename sheets: "Table 2" -> "Metabolites", "Table 3" -> "Genes"
$wb.Worksheets.Item(2).Name = "Metabolites"
$wb.Worksheets.Item(3).Name = "Genes"

$ws1 = $wb.Worksheets.Item(1)   # Table 1
$ws3 = $wb.Worksheets.Item(3)   # Genes (was Table 3)

# Update the selection on "Table 1" (sheet1) to C14, without changing
# which sheet tab is active.
$ws1.Range("C14").Select()

# "Genes" stays the active tab: update its selection to D8 and scroll the
# view down so row 4 is at the top.
$ws3.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws3.Range("D8").Select()

# Make sure Genes (the originally active tab) remains the active tab.
$ws3.Activate()
